# Insert a new row at the top of the data table (row 2) and populate it
# with the "0 months -> discount factor 1" data point, shifting all the
# existing month/discount-factor rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 2:34 down to 3:35, leaving a blank row 2.
$ws.Rows("2:2").Insert()

# Populate the newly inserted row.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1

# Match the "0.0000" number format used by the rest of column B.
$ws.Range("B2").NumberFormat = "0.0000"

# Leave the active selection on the newly added cell.
$ws.Range("B2").Select() | Out-Null
